# feat(HES-1488): Import MBUS configurations
#
# Replace the generic "MM_BUS" protocol enum value on the "Data
# specification" sheet with two dedicated values: "Mbus" and "WMBus".

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Data specification")

# K6 previously held "MM_BUS" -> rename to "Mbus"
$ws2.Range("K6").Value = "Mbus"

# K7 was empty -> add the new "WMBus" enum value
$ws2.Range("K7").Value = "WMBus"

# Reflect the user's focus moving to the "Data specification" sheet,
# with the selection left on the newly edited cell K7.
$ws2.Select()
$ws2.Range("K7").Select()
